# Edit script: apply the diff changes to the document.
#
# Change 1: In the "Extending" section's first paragraph, remove the
# mention of the AzureVM / virtual machines / AzureStor / storage accounts
# hyperlinks (and their connecting text), turning:
#   "...service-specific features; packages that do this include AzureVM
#    for virtual machines, and AzureStor for storage accounts. For example,
#    instead of using a generic ..."
# into:
#   "...service-specific features; For example, instead of using a
#    generic ..."
#
# Change 2: Remove the whole paragraph that starts with
#   "For more information, see the "Extending AzureRMR" vignette..."

$d = $word.ActiveDocument

# --- Change 1 -------------------------------------------------------

# The four hyperlinks to drop, in document order: AzureVM, virtual
# machines, AzureStor, storage accounts. Unlink them first (this removes
# the <w:hyperlink> wrapper but keeps the display text in place), then
# delete the now-plain text together with the connective wording around
# it in a single pass so the surrounding runs merge back together.
# Unlinking from the last hyperlink back to the first keeps the document
# tidy (no leftover empty spell-check markers).
$linkText = @("AzureVM", "virtual machines", "AzureStor", "storage accounts")
$toUnlink = @()
for ($i = 1; $i -le $d.Hyperlinks.Count; $i++) {
    $h = $d.Hyperlinks.Item($i)
    if ($linkText -contains $h.Range.Text) {
        $toUnlink += $i
    }
}
for ($i = $toUnlink.Count - 1; $i -ge 0; $i--) {
    $d.Hyperlinks.Item($toUnlink[$i]).Delete()
}

$startMarker = $d.Content
$startMarker.Find.Execute("packages that do this include ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$startPos = $startMarker.Start

$endMarker = $d.Content
$endMarker.Find.Execute("For example, instead of using a generic ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$endPos = $endMarker.Start

$d.Range($startPos, $endPos).Delete()

# --- Change 2 -------------------------------------------------------

# Remove the entire "For more information, see the..." paragraph,
# including its hyperlink to the AzureRMR vignette.
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*For more information, see the*AzureRMR*vignette*") {
        $d.Range($p.Range.Start, $p.Range.End).Delete()
        break
    }
}
